{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the six small text corrections described by the diff:\n//  1. \"Groups:\" paragraph \u2014 rewrite the trailing sentence.\n//  2. Add a missing period after \"...DC power and AC power\"\n//  3. Insert \" (one column)\" before the \"pandas DataFrame\" mention.\n//  4. \"...For the first 4 tasks:\" -> \"...tasks, the following points will be assessed:\"\n//  5. \"For the 5th task: Correctness of the calculation\" -> \"...task: The correctness of the calculation is assessed.\"\n//  6. \"...15% for each failure\" -> \"...failure, the following points will be assessed:\"\n\nasync function findUnique(body, text) {\n  const results = body.search(text, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for ${JSON.stringify(text)}, found ${results.items.length}`);\n  }\n  return results.items[0];\n}\n\nconst body = context.document.body;\n\n// 1. \"Groups:\" paragraph typo fix.\nlet range = await findUnique(body, \" Make groups of 2 people, you can also work alone.\");\nrange.insertText(\n  \" Groups of 2 people, working alone is also an option (and highly recommended).\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2. Add the missing period at the end of the \"Expected outputs\" sentence.\nrange = await findUnique(body, \"POA irradiance, module temperature, DC power and AC power\");\nrange.insertText(\n  \"POA irradiance, module temperature, DC power and AC power.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 3. Mention that the DataFrame has one column.\nrange = await findUnique(body, \"which contains a pandas DataFrame\");\nrange.insertText(\n  \"which contains a (one column) pandas DataFrame\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 4. \"For the first 4 tasks:\" -> add the \"the following points will be assessed\" clause.\nrange = await findUnique(body, \" tasks:\");\nrange.insertText(\" tasks, the following points will be assessed:\", \"Replace\");\nawait context.sync();\n\n// 5. \"For the 5th task: Correctness of the calculation\" -> full sentence.\nrange = await findUnique(body, \" task: Correctness of the calculation\");\nrange.insertText(\" task: The correctness of the calculation is assessed.\", \"Replace\");\nawait context.sync();\n\n// 6. \"15% for each failure\" -> add the \"the following points will be assessed:\" clause.\nrange = await findUnique(body, \"% for each failure\");\nrange.insertText(\"% for each failure, the following points will be assessed:\", \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the six small text corrections described by the diff using\n# Find & Replace over the whole document story, same as a human doing\n# Ctrl+H for each typo fix.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $result = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1. \"Groups:\" paragraph typo fix.\nReplace-Text \" Make groups of 2 people, you can also work alone.\" \" Groups of 2 people, working alone is also an option (and highly recommended).\"\n\n# 2. Add the missing period at the end of the \"Expected outputs\" sentence.\nReplace-Text \"POA irradiance, module temperature, DC power and AC power\" \"POA irradiance, module temperature, DC power and AC power.\"\n\n# 3. Mention that the DataFrame has one column.\nReplace-Text \"which contains a pandas DataFrame\" \"which contains a (one column) pandas DataFrame\"\n\n# 4. \"For the first 4 tasks:\" -> add the \"the following points will be assessed\" clause.\nReplace-Text \" tasks:\" \" tasks, the following points will be assessed:\"\n\n# 5. \"For the 5th task: Correctness of the calculation\" -> full sentence.\nReplace-Text \" task: Correctness of the calculation\" \" task: The correctness of the calculation is assessed.\"\n\n# 6. \"15% for each failure\" -> add the \"the following points will be assessed:\" clause.\nReplace-Text \"% for each failure\" \"% for each failure, the following points will be assessed:\"\n"}
